$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Update the cached "datetimeFigureOut" footer date text from
#    4/3/22 to 4/5/22 on the slide master and on every slide layout.
# ---------------------------------------------------------------------
function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "4/3/22") {
                $tr.Text = "4/5/22"
            }
        }
    }
}

Update-DatePlaceholder $p.SlideMaster.Shapes

$layouts = $p.SlideMaster.CustomLayouts
for ($L = 1; $L -le $layouts.Count; $L++) {
    Update-DatePlaceholder $layouts.Item($L).Shapes
}

# ---------------------------------------------------------------------
# 2) Recolor the three rounded rectangles ("Rectangle 163",
#    "Rectangle 86", "Rectangle 95") on slides 1 and 2 from
#    C9146F to B58900. (The same-colored "Oval 6" on slide 4 stays
#    untouched.)
# ---------------------------------------------------------------------
$targetNames = @("Rectangle 163", "Rectangle 86", "Rectangle 95")

foreach ($slideIndex in 1, 2) {
    $slide = $p.Slides.Item($slideIndex)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($targetNames -contains $shp.Name) {
            $shp.Fill.ForeColor.RGB = 0x0089B5
        }
    }
}
